# CIV-6625 Update GA order template
# Removes the "Classification: Controlled" watermark text box that was
# anchored in the default (primary) footer of the document.

$d = $word.ActiveDocument

$wdHeaderFooterPrimary = 1

foreach ($sec in $d.Sections) {
    $ftr = $sec.Footers.Item($wdHeaderFooterPrimary)
    if ($ftr.Exists) {
        for ($i = $ftr.Shapes.Count; $i -ge 1; $i--) {
            $shape = $ftr.Shapes.Item($i)
            $shape.Delete()
        }
    }
}
